$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the Heading1 title
#    paragraph ("Play Cornelius Free: Fun Features & Cat Characters").
#    The new paragraph needs: a leading empty run, a bold "Meta description"
#    run, and a plain run with the rest of the sentence.
# ---------------------------------------------------------------------------

$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item(2)
$newPara.Style = "Normal"

# Build the text (with bold split) in a scratch paragraph right after the new
# (still-empty) paragraph, so that when we copy/paste it back into the target
# paragraph, the target's own pre-existing empty run ends up first (matching
# the leading <w:r/> pattern used throughout this document).
$newPara.Range.InsertParagraphAfter()
$scratchPara = $d.Paragraphs.Item(3)
$scratchPara.Style = "Normal"
$scratchPara.Range.Text = "Meta description: Read our Cornelius review and play for free. Enjoy medium volatility, special features, fun cat characters, and appetizing instant prizes in this game from NetEnt."

$scratchPara = $d.Paragraphs.Item(3)
$boldStart = $scratchPara.Range.Start
$boldEnd = $boldStart + ("Meta description").Length
$boldRange = $d.Range($boldStart, $boldEnd)
$boldRange.Bold = 1

$scratchPara = $d.Paragraphs.Item(3)
$copyRange = $d.Range($scratchPara.Range.Start, $scratchPara.Range.End - 1)
$copyRange.Copy()

$targetPara = $d.Paragraphs.Item(2)
$targetPara.Range.Paste()

# Remove the scratch paragraph now that its content has been copied over.
$scratchPara = $d.Paragraphs.Item(3)
$scratchPara.Range.Delete()

# ---------------------------------------------------------------------------
# 2) At the end of the document: drop the duplicated "Play Cornelius Free:
#    Fun Features & Cat Characters" (bold) paragraph entirely, and rewrite
#    the italic paragraph's text with the new image-prompt copy.
# ---------------------------------------------------------------------------

$count = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs.Item($count - 1)
$dupTitlePara.Range.Delete()

$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastTextRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
$lastTextRange.Text = "Create a cartoon-style feature image showcasing Cornelius, the gluttonous kitty with glasses. The image should feature a happy Maya warrior with glasses, who is enjoying a sweet treat with Cornelius. The Maya warrior should be dressed in traditional warrior attire with a headdress, while Cornelius should be sitting on a giant cookie, wearing his green shirt that's too tight for him. The background should be a mint green kitchen with a fridge on the right and Cornelius' red bowl on the left. The image should have bright and cheerful colors to reflect the fun theme of the game."
